# Version_13_01 Added code to handle functions
#
# 1) Rename "Bank Manager Login" -> "BankManagerLogin"
# 2) Insert a new "DeleteCustomer" sheet between "BankManagerLogin" and "AddCustomer"
#    populated with the txtSearchCustomer object-repository row
# 3) Fix up the selection on the BankManagerLogin sheet
# 4) Adjust the view on the OpenAccount sheet (scrolled, new active cell)

$wb = $excel.ActiveWorkbook

# --- 1) rename Bank Manager Login sheet -------------------------------------------------
$wsManager = $wb.Worksheets.Item("Bank Manager Login")
$wsManager.Name = "BankManagerLogin"

# --- 2) insert the new DeleteCustomer sheet before AddCustomer --------------------------
$wsAddCustomer = $wb.Worksheets.Item("AddCustomer")
$wsDelete = $wb.Worksheets.Add($wsAddCustomer)
$wsDelete.Name = "DeleteCustomer"

# Header row (re-uses the shared "Object"/"ObjectID"/"Locator"/"Type" strings)
$wsDelete.Range("A1").Value = "Object"
$wsDelete.Range("B1").Value = "ObjectID"
$wsDelete.Range("C1").Value = "Locator"
$wsDelete.Range("D1").Value = "Type"
$wsDelete.Range("A1:D1").Interior.Color = 49407

# Data row
$wsDelete.Range("A2").Value = "txtSearchCustomer"
$wsDelete.Range("B2").Value = 1
$wsDelete.Range("C2").Value = "'//input[@placeholder='Search Customer']"
$wsDelete.Range("D2").Value = "TextBox"

$wsDelete.Columns.Item(3).ColumnWidth = 35.5

[void]$wsDelete.Range("F9").Select()

# --- 3) fix selection on BankManagerLogin ------------------------------------------------
[void]$wsManager.Range("A1:D1").Select()

# --- 4) adjust view on OpenAccount sheet -------------------------------------------------
$wsOpenAccount = $wb.Worksheets.Item("OpenAccount")
$wsOpenAccount.Activate()
[void]$wsOpenAccount.Range("C6").Select()

# --- leave DeleteCustomer as the active sheet/tab at the end -----------------------------
$wsDelete.Activate()
